$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 185.66667
$ws.Range("I6").Value = 171.375
$ws.Range("K6").Value = 514.125
$ws.Range("M6").Value = -402.125
$ws.Range("H12").Value = 725
$ws.Range("H15").Value = 822.1818
$ws.Range("I15").Value = 822.1818
$ws.Range("K15").Value = 2466.5454
$ws.Range("M15").Value = -2297.5454
$ws.Range("H21").Value = 37586
$ws.Range("I21").Value = 37586
$ws.Range("K21").Value = 37586
$ws.Range("M21").Value = -37118
$ws.Range("H23").Value = 37586
$ws.Range("I23").Value = 37586
$ws.Range("K23").Value = 37586
$ws.Range("M23").Value = -37352
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 33469.5
$ws.Range("I38").Value = 33469.5
$ws.Range("K38").Value = 100408.5
$ws.Range("M38").Value = -100036.5
$ws.Range("H58").Value = 2734.125
$ws.Range("J58").Value = 6666.6665
$ws.Range("L58").Value = 19999.9995
$ws.Range("N58").Value = -20299.9995
$ws.Range("H86").Value = 5160.95
$ws.Range("J86").Value = 8609
$ws.Range("L86").Value = 8609
$ws.Range("N86").Value = -10855
$ws.Range("H89").Value = 5160.95
$ws.Range("J89").Value = 8609
$ws.Range("L89").Value = 43045
$ws.Range("N89").Value = -54277
$ws.Range("H94").Value = 2653.5
$ws.Range("I94").Value = 2653.5
$ws.Range("K94").Value = 2653.5
$ws.Range("M94").Value = -2202.5
$ws.Range("H112").Value = 3095.1428
$ws.Range("J112").Value = 3854.5
$ws.Range("L112").Value = 11563.5
$ws.Range("N112").Value = -13779.5
$ws.Range("H137").Value = 2238.3333
$ws.Range("I137").Value = 2234
$ws.Range("K137").Value = 6702
$ws.Range("M137").Value = -4152
$ws.Range("H138").Value = 7554.7827
$ws.Range("I138").Value = 3943.6924
$ws.Range("J138").Value = 12249.2
$ws.Range("K138").Value = 11831.0772
$ws.Range("L138").Value = 36747.60000000001
$ws.Range("M138").Value = -6691.0772
$ws.Range("N138").Value = -47027.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 20910494
$ws.Range("I61").Value = 35001348
$ws.Range("K61").Value = 35001348
$ws.Range("M61").Value = -35001136
$ws.Range("H136").Value = 20910494
$ws.Range("I136").Value = 35001348
$ws.Range("K136").Value = 105004044
$ws.Range("M136").Value = -105001494

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 9999999
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 43481636
$ws.Range("I31").Value = 76926740
$ws.Range("J31").Value = 2999.5
$ws.Range("K31").Value = 76926740
$ws.Range("L31").Value = 2999.5
$ws.Range("M31").Value = -76926445
$ws.Range("N31").Value = -3589.5
$ws.Range("H34").Value = 43481636
$ws.Range("I34").Value = 76926740
$ws.Range("J34").Value = 2999.5
$ws.Range("K34").Value = 76926740
$ws.Range("L34").Value = 2999.5
$ws.Range("M34").Value = -76926538
$ws.Range("N34").Value = -3403.5
$ws.Range("H58").Value = 2989.8667
$ws.Range("I58").Value = 3030.1667
$ws.Range("K58").Value = 3030.1667
$ws.Range("M58").Value = -2827.1667
$ws.Range("H94").Value = 942.7059
$ws.Range("I94").Value = 847.875
$ws.Range("K94").Value = 847.875
$ws.Range("M94").Value = -396.875
$ws.Range("H107").Value = 1291.5518
$ws.Range("I107").Value = 1048.8334
$ws.Range("K107").Value = 1048.8334
$ws.Range("M107").Value = 871.1666
$ws.Range("H136").Value = 2989.8667
$ws.Range("I136").Value = 3030.1667
$ws.Range("K136").Value = 9090.500100000001
$ws.Range("M136").Value = -6540.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 72.125
$ws.Range("J2").Value = 77.666664
$ws.Range("L2").Value = 465.999984
$ws.Range("N2").Value = -691.999984
$ws.Range("H5").Value = 2601.6667
$ws.Range("J5").Value = 2354
$ws.Range("L5").Value = 7062
$ws.Range("N5").Value = -7286
$ws.Range("H17").Value = 11532.5
$ws.Range("J17").Value = 11532.5
$ws.Range("L17").Value = 34597.5
$ws.Range("N17").Value = -34935.5
$ws.Range("H34").Value = 3259.5833
$ws.Range("I34").Value = 210.4
$ws.Range("J34").Value = 5437.5713
$ws.Range("K34").Value = 631.2
$ws.Range("L34").Value = 16312.7139
$ws.Range("M34").Value = -547.2
$ws.Range("N34").Value = -16480.7139
$ws.Range("H39").Value = 10838.833
$ws.Range("J39").Value = 12606.6
$ws.Range("L39").Value = 37819.8
$ws.Range("N39").Value = -38407.8
$ws.Range("H55").Value = 5155.125
$ws.Range("J55").Value = 8583.143
$ws.Range("L55").Value = 25749.429
$ws.Range("N55").Value = -26103.429
$ws.Range("H68").Value = 2618.2778
$ws.Range("I68").Value = 2554.2
$ws.Range("J68").Value = 2698.375
$ws.Range("K68").Value = 7662.599999999999
$ws.Range("L68").Value = 8095.125
$ws.Range("M68").Value = -6851.599999999999
$ws.Range("N68").Value = -9717.125
$ws.Range("H71").Value = 2618.2778
$ws.Range("I71").Value = 2554.2
$ws.Range("J71").Value = 2698.375
$ws.Range("K71").Value = 22987.8
$ws.Range("L71").Value = 24285.375
$ws.Range("M71").Value = -18931.8
$ws.Range("N71").Value = -32397.375
$ws.Range("H122").Value = 167302
$ws.Range("I122").Value = 167302
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1505718
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1503268
$ws.Range("N122").ClearContents()
$ws.Range("H135").Value = 2601.6667
$ws.Range("J135").Value = 2354
$ws.Range("L135").Value = 21186
$ws.Range("N135").Value = -26256
$ws.Range("H137").Value = 4334.0586
$ws.Range("I137").Value = 1647.909
$ws.Range("K137").Value = 4943.727000000001
$ws.Range("M137").Value = 156.2729999999992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3575228
$ws.Range("I132").Value = 3537.2856
$ws.Range("J132").Value = 14290300
$ws.Range("K132").Value = 10611.8568
$ws.Range("L132").Value = 42870900
$ws.Range("M132").Value = -8081.856800000001
$ws.Range("N132").Value = -42875960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6254.5454
$ws.Range("I16").Value = 2980
$ws.Range("J16").Value = 10184
$ws.Range("K16").Value = 2980
$ws.Range("L16").Value = 10184
$ws.Range("M16").Value = -2810
$ws.Range("N16").Value = -10524
$ws.Range("H22").Value = 1300.5
$ws.Range("I22").Value = 1101
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 1101
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -806
$ws.Range("N22").Value = -2090
$ws.Range("H27").Value = 1300.5
$ws.Range("I27").Value = 1101
$ws.Range("J27").Value = 1500
$ws.Range("K27").Value = 1101
$ws.Range("L27").Value = 1500
$ws.Range("M27").Value = -994
$ws.Range("N27").Value = -1714
$ws.Range("H132").Value = 3888.3157
$ws.Range("I132").Value = 2621.6155
$ws.Range("K132").Value = 7864.8465
$ws.Range("M132").Value = -5334.8465
$ws.Range("H136").Value = 4924.2
$ws.Range("I136").Value = 5088.1113
$ws.Range("J136").Value = 3449
$ws.Range("K136").Value = 15264.3339
$ws.Range("L136").Value = 10347
$ws.Range("M136").Value = -12714.3339
$ws.Range("N136").Value = -15447

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18422.4
$ws.Range("J62").Value = 20978
$ws.Range("L62").Value = 20978
$ws.Range("N62").Value = -22226
$ws.Range("H65").Value = 18422.4
$ws.Range("J65").Value = 20978
$ws.Range("L65").Value = 104890
$ws.Range("N65").Value = -111130
$ws.Range("H126").Value = 2519.5625
$ws.Range("J126").Value = 1354
$ws.Range("L126").Value = 4062
$ws.Range("N126").Value = -9002
